# Auto-generated edit script: update crypto price/volume columns (D, E)
# to refreshed values per the Thu Jan  4 09:51:03 UTC 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.049.11"
$ws.Range("E2").Value = "  -5.41%  "
$ws.Range("D3").Value = "2.225.65"
$ws.Range("E3").Value = "  -6.45%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'322.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("D6").Value = "'98.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.77%  "
$ws.Range("E7").Value = "  -8.89%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.87%  "
$ws.Range("D10").Value = "'36.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.65%  "
$ws.Range("D11").Value = "'54.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("D12").Value = "'0.0831"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.84%  "
$ws.Range("E13").Value = "  -10.73%  "
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").Value = "2.566.10"
$ws.Range("E15").Value = "  -6.51%  "
$ws.Range("E16").Value = "  -12.42%  "
$ws.Range("D17").Value = "'14.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.45%  "
$ws.Range("D18").Value = "2.226.72"
$ws.Range("E18").Value = "  -6.51%  "
$ws.Range("D19").Value = "42.980.73"
$ws.Range("D20").Value = "'13.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.44%  "
$ws.Range("D21").Value = "0.0₃0966"
$ws.Range("E21").Value = "  -9.57%  "
$ws.Range("E22").Value = "  -11.00%  "
$ws.Range("D23").Value = "'3.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.55%  "
$ws.Range("D24").Value = "'65.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.18%  "
$ws.Range("D25").Value = "'236.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.61%  "
$ws.Range("D26").Value = "'2.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.32%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").Value = "'10.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.63%  "
$ws.Range("D31").Value = "'6.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.84%  "
$ws.Range("D32").Value = "'36.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("D33").Value = "'20.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.74%  "
$ws.Range("E34").Value = "  -9.33%  "
$ws.Range("D35").Value = "'155.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.50%  "
$ws.Range("E36").Value = "  -7.47%  "
$ws.Range("D37").Value = "'3.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'0.121"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.03%  "
$ws.Range("D39").Value = "'1.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("E40").Value = "  -7.45%  "
$ws.Range("E41").Value = "  -11.60%  "
$ws.Range("E42").Value = "  -8.20%  "
$ws.Range("E43").Value = "  -9.75%  "
$ws.Range("D44").Value = "'14.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.93%  "
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "1.734.74"
$ws.Range("E46").Value = "  -7.35%  "
$ws.Range("D47").Value = "'84.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -14.01%  "
$ws.Range("E48").Value = "  -12.09%  "
$ws.Range("D49").Value = "'8.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("E50").Value = "  -13.52%  "
$ws.Range("D51").Value = "'74.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -12.83%  "
